$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for account 005654767 / DIEGO / 1265.69 (row 27),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(27).Delete()

# Update the balance for account 004983378 / MARCELO, which was on row 150
# and is now row 149 after the deletion above, from 99.59 to 99.33.
$ws.Cells.Item(149, 3).Value = 99.33
